$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case schema ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case municipality / state names (mirrors upstream cleaning script) ---
$ws.Range("B17").Value = "Amatenango De La Frontera"
$ws.Range("B20").Value = "Bejucal De Ocampo"
$ws.Range("B25").Value = "Comitán De Domínguez"
$ws.Range("B36").Value = "Mazapa De Madero"
$ws.Range("B42").Value = "San Cristóbal De Las Casas"
$ws.Range("B58").Value = "Hidalgo Del Parral"
$ws.Range("B64").Value = "San Juan De Sabinas"
$ws.Range("A67").Value = "Ciudad De México"
$ws.Range("B71").Value = "Cuajimalpa De Morelos"
$ws.Range("B86").Value = "Pánuco De Coronado"
$ws.Range("A91").Value = "Estado De México"
$ws.Range("B91").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B92").Value = "Almoloya De Alquisiras"
$ws.Range("B93").Value = "Almoloya De Juárez"
$ws.Range("B98").Value = "Atizapán De Zaragoza"
$ws.Range("B107").Value = "Ecatepec De Morelos"
$ws.Range("B112").Value = "Ixtapan De La Sal"
$ws.Range("B121").Value = "Naucalpan De Juárez"
$ws.Range("B127").Value = "San Felipe Del Progreso"
$ws.Range("B128").Value = "San Simón De Guerero"
$ws.Range("B129").Value = "Soyaniquilpan De Juárez"
$ws.Range("B135").Value = "Tenango Del Valle"
$ws.Range("B142").Value = "Tlalnepantla De Baz"
$ws.Range("B147").Value = "Valle De Bravo"
$ws.Range("B148").Value = "Valle De Chalco Solidaridad"
$ws.Range("B149").Value = "Villa De Allende"
$ws.Range("B150").Value = "Villa Del Carbón"
$ws.Range("A158").Value = "Guanajuato"
$ws.Range("B160").Value = "Apaseo El Grande"
$ws.Range("B167").Value = "Jaral Del Progreso"
$ws.Range("B177").Value = "San Francisco Del Rincón"
$ws.Range("B179").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B181").Value = "Silao De La Victoria"
$ws.Range("B185").Value = "Valle De Santiago"
$ws.Range("B190").Value = "Acapulco De Juárez"
$ws.Range("B192").Value = "Ajuchitlán Del Progreso"
$ws.Range("B196").Value = "Atenango Del Río"
$ws.Range("B197").Value = "Atoyac De Álvarez"
$ws.Range("B198").Value = "Ayutla De Los Libres"
$ws.Range("B201").Value = "Chilapa De Álvarez"
$ws.Range("B202").Value = "Chilpancingo De Los Bravo"
$ws.Range("B203").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B206").Value = "Coyuca De Benítez"
$ws.Range("B207").Value = "Coyuca De Catalán"
$ws.Range("B211").Value = "Cuetzala Del Progreso"
$ws.Range("B212").Value = "Cutzamala De Pinzón"
$ws.Range("B217").Value = "Huitzuco De Los Figueroa"
$ws.Range("B218").Value = "Iguala De La Independencia"
$ws.Range("B220").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B221").Value = "Zihuatanejo De Azueta"
$ws.Range("B233").Value = "Taxco De Alarcón"
$ws.Range("B235").Value = "Técpan De Galeana"
$ws.Range("B237").Value = "Tepecoacuilco De Trujano"
$ws.Range("B239").Value = "Tlapa De Comonfort"
$ws.Range("B252").Value = "Atotonilco El Grande"
$ws.Range("B255").Value = "Cuautepec De Hinojosa"
$ws.Range("B257").Value = "Huejutla De Reyes"
$ws.Range("B261").Value = "Mineral Del Monte"
$ws.Range("B263").Value = "Omitlán De Juárez"
$ws.Range("B264").Value = "Pachuca De Soto"
$ws.Range("B266").Value = "Progreso De Obregón"
$ws.Range("B267").Value = "Santiago Tulantepec De Lugo Guerero"
$ws.Range("B270").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B272").Value = "Tezontepec De Aldama"
$ws.Range("B275").Value = "Tulancingo De Bravo"
$ws.Range("B276").Value = "Zacualtipán De Ángeles"
$ws.Range("B277").Value = "Zapotlán De Juárez"
$ws.Range("B280").Value = "Ahualulco De Mercado"
$ws.Range("B281").Value = "Atotonilco El Alto"
$ws.Range("B282").Value = "Autlán De Navarro"
$ws.Range("B286").Value = "Encarnación De Díaz"
$ws.Range("B288").Value = "Huejuquilla El Alto"
$ws.Range("B291").Value = "Lagos De Moreno"
$ws.Range("B309").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B355").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B359").Value = "Tetela Del Volcán"
$ws.Range("B364").Value = "Ixtlán Del Río"
$ws.Range("B378").Value = "Constancia Del Rosario"
$ws.Range("B380").Value = "Fresnillo De Trujano"
$ws.Range("B381").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B382").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B383").Value = "Ixtlán De Juárez"
$ws.Range("B384").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B386").Value = "Mariscala De Juárez"
$ws.Range("B388").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B389").Value = "Oaxaca De Juárez"
$ws.Range("B390").Value = "Ocotlán De Morelos"
$ws.Range("B391").Value = "Putla Villa De Guerero"
$ws.Range("B392").Value = "Reforma De Pineda"
$ws.Range("B403").Value = "San Dionisio Del Mar"
$ws.Range("B405").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B408").Value = "San Juan De Los Cués"
$ws.Range("B416").Value = "San Miguel Del Puerto"
$ws.Range("B445").Value = "Santo Domingo De Morelos"
$ws.Range("B449").Value = "Teotitlán De Flores Magón"
$ws.Range("B450").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B451").Value = "Tlacolula De Matamoros"
$ws.Range("B452").Value = "Tlalixtac De Cabrera"
$ws.Range("B453").Value = "Totontepec Villa De Morelos"
$ws.Range("B454").Value = "Villa De Chilapa De Díaz"
$ws.Range("B455").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B468").Value = "Ayotoxco De Guerero"
$ws.Range("B471").Value = "Chalchicomula De Sesma"
$ws.Range("B492").Value = "Huehuetlán El Chico"
$ws.Range("B495").Value = "Huitzilan De Serdán"
$ws.Range("B496").Value = "Ixcamilpa De Guerero"
$ws.Range("B498").Value = "Izúcar De Matamoros"
$ws.Range("B504").Value = "Los Reyes De Juárez"
$ws.Range("B508").Value = "Palmar De Bravo"
$ws.Range("B519").Value = "San Nicolás De Los Ranchos"
$ws.Range("B522").Value = "San Salvador El Verde"
$ws.Range("B527").Value = "Tecali De Herrera"
$ws.Range("B534").Value = "Tepango De Rodríguez"
$ws.Range("B538").Value = "Tepexi De Rodríguez"
$ws.Range("B540").Value = "Teteles De Avila Castillo"
$ws.Range("B564").Value = "Cadereyta De Montes"
$ws.Range("B565").Value = "Jalpan De Serra"
$ws.Range("B566").Value = "Pinal De Amoles"
$ws.Range("B572").Value = "Ciudad Del Maíz"
$ws.Range("B604").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B617").Value = "Tetla De La Solidaridad"
$ws.Range("B626").Value = "Boca Del Río"
$ws.Range("B635").Value = "Cosamaloapan De Carpio"
$ws.Range("B639").Value = "Ixhuatlán Del Café"
$ws.Range("B647").Value = "Martínez De La Torre"
$ws.Range("B651").Value = "Mixtla De Altamirano"
$ws.Range("B655").Value = "Paso De Ovejas"
$ws.Range("B659").Value = "Poza Rica De Hidalgo"
$ws.Range("B663").Value = "Soledad De Doblado"
$ws.Range("B679").Value = "Zozocolco De Hidalgo"
$ws.Range("B687").Value = "Jiménez Del Teul"
$ws.Range("B690").Value = "Teúl De González Ortega"
$ws.Range("B691").Value = "Villa De Cos"

# --- Tiny floating point recompute on D230 ---
$ws.Range("D230").Value = 0.009308098045299413

# --- Drop trailing footnote/metadata rows (696:700) ---
$ws.Range("A696:D700").EntireRow.Delete()

